$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 57.319636
$ws.Range("H2").Value = 171.958908
$ws.Range("I2").Value = 0.5476981520382651
$ws.Range("J2").Value = 0.5476981520382651
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 216.578748128204
$ws.Range("R2").Value = 1949.208733153836
$ws.Range("S2").Value = 0.2329344646734626
$ws.Range("T2").Value = 0.2329344646734626

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.319636
$ws.Range("H3").Value = 171.958908
$ws.Range("I3").Value = 0.5476981520382651
$ws.Range("J3").Value = 0.5476981520382651
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 248.3891017078533
$ws.Range("R3").Value = 2235.50191537068
$ws.Range("S3").Value = 0.2671470905482921
$ws.Range("T3").Value = 0.2671470905482921

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 57.319636
$ws.Range("H4").Value = 171.958908
$ws.Range("I4").Value = 0.5476981520382651
$ws.Range("J4").Value = 0.5476981520382651
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 16.67917338800533
$ws.Range("R4").Value = 150.112560492048
$ws.Range("S4").Value = 0.01793876064899528
$ws.Range("T4").Value = 0.01793876064899528

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.319636
$ws.Range("H5").Value = 171.958908
$ws.Range("I5").Value = 0.5476981520382651
$ws.Range("J5").Value = 0.5476981520382651
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 27.59397847512533
$ws.Range("R5").Value = 248.345806276128
$ws.Range("S5").Value = 0.02967783616751516
$ws.Range("T5").Value = 0.02967783616751516

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.975184333333334
$ws.Range("H6").Value = 17.925553
$ws.Range("I6").Value = 0.05709382762749331
$ws.Range("J6").Value = 0.05709382762749331
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 22.57686951725567
$ws.Range("R6").Value = 203.191825655301
$ws.Range("S6").Value = 0.02428184233427896
$ws.Range("T6").Value = 0.02428184233427896

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.975184333333334
$ws.Range("H7").Value = 17.925553
$ws.Range("I7").Value = 0.05709382762749331
$ws.Range("J7").Value = 0.05709382762749331
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 25.89288370734778
$ws.Range("R7").Value = 233.03595336613
$ws.Range("S7").Value = 0.02784827716176942
$ws.Range("T7").Value = 0.02784827716176943

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.975184333333334
$ws.Range("H8").Value = 17.925553
$ws.Range("I8").Value = 0.05709382762749331
$ws.Range("J8").Value = 0.05709382762749331
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 1.738691004963111
$ws.Range("R8").Value = 15.648219044668
$ws.Range("S8").Value = 0.00186999445686105
$ws.Range("T8").Value = 0.00186999445686105

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.975184333333334
$ws.Range("H9").Value = 17.925553
$ws.Range("I9").Value = 0.05709382762749331
$ws.Range("J9").Value = 0.05709382762749331
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 2.876485605716444
$ws.Range("R9").Value = 25.888370451448
$ws.Range("S9").Value = 0.003093713674583871
$ws.Range("T9").Value = 0.003093713674583872

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.32302533333333
$ws.Range("H10").Value = 96.969076
$ws.Range("I10").Value = 0.3088515991858827
$ws.Range("J10").Value = 0.3088515991858827
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 122.1305795174547
$ws.Range("R10").Value = 1099.175215657092
$ws.Range("S10").Value = 0.1313537057815016
$ws.Range("T10").Value = 0.1313537057815016

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 32.32302533333333
$ws.Range("H11").Value = 96.969076
$ws.Range("I11").Value = 0.3088515991858827
$ws.Range("J11").Value = 0.3088515991858827
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 140.0687057228844
$ws.Range("R11").Value = 1260.61835150596
$ws.Range("S11").Value = 0.1506464935597068
$ws.Range("T11").Value = 0.1506464935597068

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 32.32302533333333
$ws.Range("H12").Value = 96.969076
$ws.Range("I12").Value = 0.3088515991858827
$ws.Range("J12").Value = 0.3088515991858827
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 9.40552630096178
$ws.Range("R12").Value = 84.64973670865601
$ws.Range("S12").Value = 0.01011581816231487
$ws.Range("T12").Value = 0.01011581816231487

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 32.32302533333333
$ws.Range("H13").Value = 96.969076
$ws.Range("I13").Value = 0.3088515991858827
$ws.Range("J13").Value = 0.3088515991858827
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 15.56047678493511
$ws.Range("R13").Value = 140.044291064416
$ws.Range("S13").Value = 0.0167355816823594
$ws.Range("T13").Value = 0.0167355816823594

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.037676333333334
$ws.Range("H14").Value = 27.113029
$ws.Range("I14").Value = 0.08635642114835883
$ws.Range("J14").Value = 0.08635642114835884
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 34.14830872724367
$ws.Range("R14").Value = 307.334778545193
$ws.Range("S14").Value = 0.03672714004319605
$ws.Range("T14").Value = 0.03672714004319606

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.037676333333334
$ws.Range("H15").Value = 27.113029
$ws.Range("I15").Value = 0.08635642114835883
$ws.Range("J15").Value = 0.08635642114835884
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 39.16389674845444
$ws.Range("R15").Value = 352.47507073609
$ws.Range("S15").Value = 0.04212149808081747
$ws.Range("T15").Value = 0.04212149808081748

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.037676333333334
$ws.Range("H16").Value = 27.113029
$ws.Range("I16").Value = 0.08635642114835883
$ws.Range("J16").Value = 0.08635642114835884
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 2.629831260413778
$ws.Range("R16").Value = 23.668481343724
$ws.Range("S16").Value = 0.00282843234675733
$ws.Range("T16").Value = 0.00282843234675733

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.037676333333334
$ws.Range("H17").Value = 27.113029
$ws.Range("I17").Value = 0.08635642114835883
$ws.Range("J17").Value = 0.08635642114835884
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 4.350785587807111
$ws.Range("R17").Value = 39.157070290264
$ws.Range("S17").Value = 0.004679350677587969
$ws.Range("T17").Value = 0.004679350677587971

